# Automatic update of files.
# Rotates the taxon-observation data among rows 10/11/12/13 and rows
# 16/17/18 on the "Artfynd" sheet: each row's content (Id, Taxonsorterings-
# ordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord and
# Publik kommentar) is replaced with the content that used to sit in a
# sibling row, per the recorded diff.
#
# NOTE: this runtime's PowerShell parser only binds function arguments
# positionally, so Set-RowData is called with positional args (no -Name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $A, $B, $E, $F, $G, $H, $Q, $R, $AC)

    $ws.Cells.Item($Row, 1).Value = $A    # A  - Id
    $ws.Cells.Item($Row, 2).Value = $B    # B  - Taxonsorteringsordning
    $ws.Cells.Item($Row, 5).Value = $E    # E  - TaxonId
    $ws.Cells.Item($Row, 6).Value = $F    # F  - Artnamn
    $ws.Cells.Item($Row, 7).Value = $G    # G  - Vetenskapligt namn

    if ($H -eq $null) {
        $ws.Cells.Item($Row, 8).Value = ""
    } else {
        $ws.Cells.Item($Row, 8).Value = $H   # H - Auktor
    }

    $ws.Cells.Item($Row, 17).Value = $Q   # Q  - Ost
    $ws.Cells.Item($Row, 18).Value = $R   # R  - Nord

    if ($AC -eq $null) {
        $ws.Cells.Item($Row, 29).Value = ""
    } else {
        $ws.Cells.Item($Row, 29).Value = $AC # AC - Publik kommentar
    }
}

# ---- Group 1: rows 10-13 ----
# New row 10 <= old row 13
Set-RowData 10 131187626 91828 5432 "Granticka" "Porodaedalea chrysoloma s.lat." $null 441180 7057056 $null

# New row 11 <= old row 12 (only A/Q/R actually change; rest identical already)
Set-RowData 11 131187615 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 441121 7057222 "Ringhack äldre"

# New row 12 <= old row 11 (only A/Q/R actually change; rest identical already)
Set-RowData 12 131187619 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 441118 7056975 "Ringhack äldre"

# New row 13 <= old row 10
Set-RowData 13 131187614 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 441097 7057228 "Ringhack"

# ---- Group 2: rows 16-18 ----
# New row 16 <= old row 18
Set-RowData 16 131187625 91828 5432 "Granticka" "Porodaedalea chrysoloma s.lat." $null 441119 7057196 $null

# New row 17 <= old row 16
Set-RowData 17 131187616 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 441123 7057208 "Ringhack äldre"

# New row 18 <= old row 17
Set-RowData 18 131187611 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 440913 7056941 "Ringhack färska och äldre"
